$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update row 16 to become the EDER YESID MORENO CAMPO record (the only remaining worker)
$ws.Range("C16").Value = "1003082835"
$ws.Range("D16").Value = "EDER YESID MORENO CAMPO"
$ws.Range("E16").Value = "2508"

# Delete rows 17-20 (the other 4 worker rows), shifting rows below up
$ws.Range("B17:J20").Delete()

# Update "Valor Mora" summary (row 11) to match the new totals
$ws.Range("E11").Value = 56940

# Update Cant. Trabajadores / Cant. Periodos counts (row 13)
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D is now narrower to fit the shorter remaining content
$ws.Columns("D").ColumnWidth = 27
